$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new data values for rows 911-927 (columns A, B, C)
$ws.Cells.Item(911, 1).Value = 909
$ws.Cells.Item(911, 2).Value = "ángel perdomo"
$ws.Cells.Item(911, 3).Value = 6.35
$ws.Cells.Item(912, 1).Value = 910
$ws.Cells.Item(912, 2).Value = "brent honeywell"
$ws.Cells.Item(912, 3).Value = 8.31
$ws.Cells.Item(913, 1).Value = 911
$ws.Cells.Item(913, 2).Value = "carl edwards"
$ws.Cells.Item(913, 3).Value = 11.12
$ws.Cells.Item(914, 1).Value = 912
$ws.Cells.Item(914, 2).Value = "dan camarena"
$ws.Cells.Item(914, 3).Value = 9.640000000000001
$ws.Cells.Item(915, 1).Value = 913
$ws.Cells.Item(915, 2).Value = "daniel lynch"
$ws.Cells.Item(915, 3).Value = 5.69
$ws.Cells.Item(916, 1).Value = 914
$ws.Cells.Item(916, 2).Value = "duane underwood"
$ws.Cells.Item(916, 3).Value = 4.33
$ws.Cells.Item(917, 1).Value = 915
$ws.Cells.Item(917, 2).Value = "j.b. bukauskas"
$ws.Cells.Item(917, 3).Value = 7.79
$ws.Cells.Item(918, 1).Value = 916
$ws.Cells.Item(918, 2).Value = "j.t. chargois"
$ws.Cells.Item(918, 3).Value = 2.52
$ws.Cells.Item(919, 1).Value = 917
$ws.Cells.Item(919, 2).Value = "jaime barría"
$ws.Cells.Item(919, 3).Value = 4.61
$ws.Cells.Item(920, 1).Value = 918
$ws.Cells.Item(920, 2).Value = "julio teherán"
$ws.Cells.Item(920, 3).Value = 1.8
$ws.Cells.Item(921, 1).Value = 919
$ws.Cells.Item(921, 2).Value = "lance mccullers"
$ws.Cells.Item(921, 3).Value = 3.16
$ws.Cells.Item(922, 1).Value = 920
$ws.Cells.Item(922, 2).Value = "matt boyd"
$ws.Cells.Item(922, 3).Value = 3.89
$ws.Cells.Item(923, 1).Value = 921
$ws.Cells.Item(923, 2).Value = "mike king"
$ws.Cells.Item(923, 3).Value = 3.55
$ws.Cells.Item(924, 1).Value = 922
$ws.Cells.Item(924, 2).Value = "mike wright"
$ws.Cells.Item(924, 3).Value = 5.5
$ws.Cells.Item(925, 1).Value = 923
$ws.Cells.Item(925, 2).Value = "néstor cortés"
$ws.Cells.Item(925, 3).Value = 2.9
$ws.Cells.Item(926, 1).Value = 924
$ws.Cells.Item(926, 2).Value = "travis lakins"
$ws.Cells.Item(926, 3).Value = 5.79
$ws.Cells.Item(927, 1).Value = 925
$ws.Cells.Item(927, 2).Value = "vladimir gutiérrez"
$ws.Cells.Item(927, 3).Value = 4.74

# Copy the row-label style (bold, centered, bordered) from the last existing
# data row (A910) onto the new A-column cells so formatting matches the rest
# of the table.
$ws.Range("A910").Copy()
$ws.Range("A911:A927").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
